# Edit script: rename title slide text, then duplicate it to create a new
# closing "Questions?" slide at the end of the deck.

$p = $ppt.ActivePresentation

# 1) Fix the title-slide wording: "TSQL Architecture" -> "T-SQL Architecture"
$titleSlide = $p.Slides.Item(1)
$titleSlide.Shapes.Item(1).TextFrame.TextRange.Text = "T-SQL Architecture"

# 2) Add a new final "Questions?" slide, based on the (now corrected) title
#    slide layout/formatting, moved to the end of the deck.
$slideCount = $p.Slides.Count
$dup = $titleSlide.Duplicate()
$questionsSlide = $dup.Item(1)
$questionsSlide.MoveTo($slideCount + 1)
$questionsSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"
